# Generate Report for Handback
# Refresh the handoff/handback timestamps recorded for the
# "582846e0-6450-4fe4-8443-be0225224ce4" file across the Overview,
# zh-cn and de-de sheets (row 2 in each table).

$wb = $excel.ActiveWorkbook

# --- Overview sheet --------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G2").Value = "2016-09-06 23:00:39"

# --- zh-cn sheet -------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H2").Value = "2016-09-06 23:00:34"
$zhcn.Range("K2").Value = "2016-09-06 23:00:51"

# --- de-de sheet -------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H2").Value = "2016-09-06 23:00:39"
$dede.Range("K2").Value = "2016-09-06 23:00:59"
